$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data first (previous layout used rows up to 27).
# Use ClearContents (not Clear) on the used range so column formatting
# (e.g. column C width/bestFit) is preserved.
$ws.UsedRange.ClearContents()

# Set values in the order that reproduces the target shared-string table order:
# 0=Number, 1=Is in DB?, 2=aaa, 3=bbbb, 4=a column, 5=b column
$ws.Range("C1").Value = "Number"
$ws.Range("D1").Value = "Is in DB?"
$ws.Range("A2").Value = "aaa"
$ws.Range("B2").Value = "bbbb"
$ws.Range("A1").Value = "a column"
$ws.Range("B1").Value = "b column"

# Row 2 numeric data
$ws.Range("C2").Value = 2949444438

# Rows 3-5 data
$ws.Range("C3").Value = 3949444438
$ws.Range("C4").Value = 3949444438
$ws.Range("C5").Value = 3949444438

# Rows 7-10 data (row 6 intentionally skipped)
$ws.Range("C7").Value = 3949444438
$ws.Range("C8").Value = 3949444438
$ws.Range("C9").Value = 3949444438
$ws.Range("C10").Value = 3949444438

# Selection matches target: activeCell D2, sqref D2:D10
$ws.Range("D2:D10").Select()
